$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (revert merge: C8 -> "物品", C13 -> "數量")
$ws.Range("C8").Value = "物品"
$ws.Range("C13").Value = "數量"

# Update the selected cell/active cell in the sheet view
$ws.Range("B11").Select()
